$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Insert the new "Speech invoked content access" writeup paragraph plus a
#    following line-break paragraph, right before the (almost empty) spacer
#    paragraph that currently just holds four spaces.
# ---------------------------------------------------------------------------
$anchor = $d.Paragraphs.Item(22).Range
$anchor.InsertParagraphBefore()
$anchor.InsertParagraphBefore()

$textPara = $d.Paragraphs.Item(22).Range
$textPara.Text = "    Speech invoked web content access is a suitable for people with dysfunctional hand motor-abilities. Instead of controlling the mouse or typing on the keyboard, such users are only required to say aloud the channel number or the short name of the web site they wish to visit. As a preliminary step they would have to train the system with their voice."

$brPara = $d.Paragraphs.Item(23).Range
$brPara.InsertAfter([char]11)

# ---------------------------------------------------------------------------
# 2. The paragraph that now precedes the manual page break (just before
#    "3. Plan of Implementation") gets stamped with a rendered-page-break
#    marker ahead of its existing <w:br/>. Rebuild that paragraph's OOXML
#    (same pPr / run attributes as before) with the marker added.
# ---------------------------------------------------------------------------
$pageBreakPara = $d.Paragraphs.Item(25).Range
$lrpbXml = '<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="00D551A1" w:rsidRDefault="00E7527A" w:rsidP="003A60F2"><w:pPr><w:spacing w:after="0"/><w:rPr><w:b/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:r><w:lastRenderedPageBreak/><w:br/></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$pageBreakPara.InsertXML($lrpbXml)

# ---------------------------------------------------------------------------
# 3. The Gantt-chart picture run had been tagged with the Far-East language
#    zh-CN; flip it to en-US.
# ---------------------------------------------------------------------------
$picturePara = $d.Paragraphs.Item(49).Range
$picturePara.LanguageIDFarEast = "en-US"

Write-Output "done"
